$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.249.25"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "2.063.48"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'249.14"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'0.667"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("D7").Value = "'58.88"
$ws.Range("E7").Value = "  +6.36%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.387"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").Value = "'0.0791"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").Value = "'0.110"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("D12").Value = "'15.93"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "'0.918"
$ws.Range("E13").Value = "  +16.71%  "
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "'5.84"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").Value = "2.086.52"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "'18.87"
$ws.Range("E17").Value = "  +14.67%  "
$ws.Range("D18").Value = "37.246.06"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "'75.61"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").Value = "0.0₃0914"
$ws.Range("E20").Value = "  +2.92%  "
$ws.Range("D21").Value = "'5.52"
$ws.Range("E21").Value = "  +4.41%  "
$ws.Range("D22").Value = "'239.74"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'2.47"
$ws.Range("E24").Value = "  +5.47%  "
$ws.Range("D25").Value = "'2.22"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "'9.65"
$ws.Range("E26").Value = "  +6.59%  "
$ws.Range("D27").Value = "'171.83"
$ws.Range("E27").Value = "  +2.90%  "
$ws.Range("D28").Value = "'20.31"
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("D29").Value = "'5.52"
$ws.Range("E29").Value = "  +18.63%  "
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").Value = "'1.17"
$ws.Range("E31").Value = "  +6.54%  "
$ws.Range("D32").Value = "'4.79"
$ws.Range("E32").Value = "  +9.74%  "
$ws.Range("D33").Value = "'0.0631"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("D34").Value = "'0.0885"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("E35").Value = "  +6.69%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +5.01%  "
$ws.Range("D38").Value = "'1.35"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "'5.21"
$ws.Range("E39").Value = "  +6.94%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").Value = "'3.13"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("E42").Value = "  +3.90%  "
$ws.Range("D43").Value = "'101.39"
$ws.Range("E43").Value = "  +6.77%  "
$ws.Range("E44").Value = "  +6.09%  "
$ws.Range("D45").Value = "'17.54"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("D46").Value = "'2.44"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").Value = "1.311.24"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("E48").Value = "  +17.54%  "
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").Value = "'6.96"
$ws.Range("E50").Value = "  +4.94%  "
$ws.Range("D51").Value = "2.247.87"
$ws.Range("E51").Value = "  +1.24%  "
